# Rename the "Pizza1" category to "Pizza L" on the ManageCategory sheet,
# and make that sheet the active/selected sheet (replacing AdminUsers as
# the previously active sheet), with the selection moved to E7.

$wb = $excel.ActiveWorkbook

$wsCategory = $wb.Worksheets.Item("ManageCategory")
$wsCategory.Range("A2").Value = "Pizza L"

$wsCategory.Activate()
$wsCategory.Range("E7").Select()
